$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.106.32'
$ws.Range('E2').Value = '  -2.20%  '
$ws.Range('D3').Value = '1.821.15'
$ws.Range('E3').Value = '  -1.52%  '
$ws.Range('E4').Value = '  -1.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.72'
$ws.Range('E5').Value = '  -3.02%  '
$ws.Range('E6').Value = '  -1.16%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4228'
$ws.Range('E7').Value = '  -1.99%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3666'
$ws.Range('E8').Value = '  -1.97%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07224'
$ws.Range('E9').Value = '  -2.05%  '
$ws.Range('E10').Value = '  -3.57%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.89'
$ws.Range('E11').Value = '  -3.54%  '
$ws.Range('D12').Value = '1.806.77'
$ws.Range('E12').Value = '  -2.29%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.647'
$ws.Range('E13').Value = '  -1.22%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.07080'
$ws.Range('E14').Value = '  -0.67%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.281'
$ws.Range('E15').Value = '  -3.25%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '89.55'
$ws.Range('E16').Value = '  +1.32%  '
$ws.Range('E17').Value = '  -1.21%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008833'
$ws.Range('E18').Value = '  -1.91%  '
$ws.Range('E19').Value = '  -1.14%  '
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.97'
$ws.Range('E20').Value = '  -3.21%  '
$ws.Range('B21').Value = 'WrappedBTC'
$ws.Range('C21').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D21').Value = '27.202.58'
$ws.Range('E21').Value = '  -1.87%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.109'
$ws.Range('E22').Value = '  -2.34%  '
$ws.Range('B23').Value = 'Cosmos'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.83'
$ws.Range('E23').Value = '  -2.46%  '
$ws.Range('B24').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C24').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D24').Value = '2.054.15'
$ws.Range('E24').Value = '  -1.32%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.972'
$ws.Range('E25').Value = '  -2.21%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '151.80'
$ws.Range('E26').Value = '  -2.62%  '
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.250'
$ws.Range('E27').Value = '  +4.84%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.30'
$ws.Range('E28').Value = '  -1.72%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.210'
$ws.Range('E29').Value = '  -3.53%  '
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '116.27'
$ws.Range('E30').Value = '  -2.14%  '
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08794'
$ws.Range('E31').Value = '  -1.97%  '
$ws.Range('B32').Value = 'ARBITRUM'
$ws.Range('C32').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.180'
$ws.Range('E32').Value = '  -4.28%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7411'
$ws.Range('E33').Value = '  -4.87%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.963'
$ws.Range('E34').Value = '  +1.35%  '
$ws.Range('B35').Value = 'Filecoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.422'
$ws.Range('E35').Value = '  -3.23%  '
$ws.Range('B36').Value = 'Frax'
$ws.Range('C36').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.000'
$ws.Range('E36').Value = '  -1.29%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.095'
$ws.Range('E37').Value = '  -3.73%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01963'
$ws.Range('E38').Value = '  -0.44%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05239'
$ws.Range('E39').Value = '  -2.16%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.300'
$ws.Range('E40').Value = '  +0.69%  '
$ws.Range('B41').Value = 'MXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.869'
$ws.Range('E41').Value = '  -0.50%  '
$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1691'
$ws.Range('E42').Value = '  +0.17%  '
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.5026'
$ws.Range('E43').Value = '  -2.41%  '
$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.591'
$ws.Range('E44').Value = '  -2.71%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.51'
$ws.Range('E45').Value = '  -1.65%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4748'
$ws.Range('E46').Value = '  -0.12%  '
$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '106.20'
$ws.Range('E47').Value = '  -3.21%  '
$ws.Range('B48').Value = 'PaxDollar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.000'
$ws.Range('E48').Value = '  -1.34%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06368'
$ws.Range('E49').Value = '  -1.83%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.655'
$ws.Range('E50').Value = '  -2.41%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.872'
$ws.Range('E51').Value = '  +0.97%  '
